$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: 'Bitcoin' -> 'Bitcoin'
$ws.Range('D2').Value = '67.784.44'
$ws.Range('E2').Value = '  -0.30%  '

# Row 3: 'Ethereum' -> 'Ethereum'
$ws.Range('D3').Value = '2.430.39'
$ws.Range('E3').Value = '  +0.31%  '

# Row 4: 'TetherUSD' -> 'TetherUSD'
$ws.Range('E4').Value = '  -0.09%  '

# Row 5: 'BNB' -> 'BNB'
$cell = $ws.Range('D5')
$cell.NumberFormat = "@"
$cell.Value = '553.07'
$cell.Style = "Normal"
$ws.Range('E5').Value = '  +0.47%  '

# Row 6: 'Solana' -> 'Solana'
$cell = $ws.Range('D6')
$cell.NumberFormat = "@"
$cell.Value = '160.10'
$cell.Style = "Normal"
$ws.Range('E6').Value = '  +0.58%  '

# Row 8: 'XRP' -> 'XRP'
$ws.Range('E8').Value = '  +2.72%  '

# Row 9: 'LidoStakedEther' -> 'Dogecoin'
$ws.Range('B9').Value = 'Dogecoin'
$ws.Range('C9').Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$cell = $ws.Range('D9')
$cell.NumberFormat = "@"
$cell.Value = '0.157'
$cell.Style = "Normal"
$ws.Range('E9').Value = '  +8.29%  '

# Row 10: 'Dogecoin' -> 'TRON'
$ws.Range('B10').Value = 'TRON'
$ws.Range('C10').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$cell = $ws.Range('D10')
$cell.NumberFormat = "@"
$cell.Value = '0.164'
$cell.Style = "Normal"
$ws.Range('E10').Value = '  -0.13%  '

# Row 11: 'TRON' -> 'Cardano'
$ws.Range('B11').Value = 'Cardano'
$ws.Range('C11').Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$cell = $ws.Range('D11')
$cell.NumberFormat = "@"
$cell.Value = '0.328'
$cell.Style = "Normal"
$ws.Range('E11').Value = '  -0.53%  '

# Row 12: 'Cardano' -> 'Toncoin'
$ws.Range('B12').Value = 'Toncoin'
$ws.Range('C12').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$cell = $ws.Range('D12')
$cell.NumberFormat = "@"
$cell.Value = '4.79'
$cell.Style = "Normal"
$ws.Range('E12').Value = '  +1.11%  '

# Row 13: 'Toncoin' -> 'WrappedBTC'
$ws.Range('B13').Value = 'WrappedBTC'
$ws.Range('C13').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D13').Value = '67.695.43'
$ws.Range('E13').Value = '  -0.53%  '

# Row 14: 'WrappedBTC' -> 'ShibaInu'
$ws.Range('B14').Value = 'ShibaInu'
$ws.Range('C14').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$cell = $ws.Range('D14')
$cell.NumberFormat = "@"
$cell.Value = '0.0000168'
$cell.Style = "Normal"
$ws.Range('E14').Value = '  +2.52%  '

# Row 15: 'ShibaInu' -> 'Avalanche'
$ws.Range('B15').Value = 'Avalanche'
$ws.Range('C15').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$cell = $ws.Range('D15')
$cell.NumberFormat = "@"
$cell.Value = '23.06'
$cell.Style = "Normal"
$ws.Range('E15').Value = '  +0.26%  '

# Row 16: 'Avalanche' -> 'Chainlink'
$ws.Range('B16').Value = 'Chainlink'
$ws.Range('C16').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$cell = $ws.Range('D16')
$cell.NumberFormat = "@"
$cell.Value = '10.33'
$cell.Style = "Normal"
$ws.Range('E16').Value = '  -2.22%  '

# Row 17: 'Chainlink' -> 'BitcoinCash'
$ws.Range('B17').Value = 'BitcoinCash'
$ws.Range('C17').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$cell = $ws.Range('D17')
$cell.NumberFormat = "@"
$cell.Value = '334.09'
$cell.Style = "Normal"
$ws.Range('E17').Value = '  -0.66%  '

# Row 18: 'BitcoinCash' -> 'Uniswap'
$ws.Range('B18').Value = 'Uniswap'
$ws.Range('C18').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$cell = $ws.Range('D18')
$cell.NumberFormat = "@"
$cell.Value = '6.82'
$cell.Style = "Normal"
$ws.Range('E18').Value = '  -0.76%  '

# Row 19: 'Uniswap' -> 'Polkadot'
$ws.Range('B19').Value = 'Polkadot'
$ws.Range('C19').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$cell = $ws.Range('D19')
$cell.NumberFormat = "@"
$cell.Value = '3.78'
$cell.Style = "Normal"
$ws.Range('E19').Value = '  +1.25%  '

# Row 20: 'Polkadot' -> 'Dai'
$ws.Range('B20').Value = 'Dai'
$ws.Range('C20').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$cell = $ws.Range('D20')
$cell.NumberFormat = "@"
$cell.Value = '1.00'
$cell.Style = "Normal"
$ws.Range('E20').Value = '  +0.09%  '

# Row 21: 'Dai' -> 'SuiNetwork'
$ws.Range('B21').Value = 'SuiNetwork'
$ws.Range('C21').Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$cell = $ws.Range('D21')
$cell.NumberFormat = "@"
$cell.Value = '1.86'
$cell.Style = "Normal"
$ws.Range('E21').Value = '  +2.04%  '

# Row 22: 'SuiNetwork' -> 'Litecoin'
$ws.Range('B22').Value = 'Litecoin'
$ws.Range('C22').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$cell = $ws.Range('D22')
$cell.NumberFormat = "@"
$cell.Value = '66.23'
$cell.Style = "Normal"
$ws.Range('E22').Value = '  +0.46%  '

# Row 23: 'Litecoin' -> 'NEARProtocol'
$ws.Range('B23').Value = 'NEARProtocol'
$ws.Range('C23').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$cell = $ws.Range('D23')
$cell.NumberFormat = "@"
$cell.Value = '3.62'
$cell.Style = "Normal"
$ws.Range('E23').Value = '  +1.37%  '

# Row 24: 'NEARProtocol' -> 'Aptos'
$ws.Range('B24').Value = 'Aptos'
$ws.Range('C24').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$cell = $ws.Range('D24')
$cell.NumberFormat = "@"
$cell.Value = '8.09'
$cell.Style = "Normal"
$ws.Range('E24').Value = '  +1.97%  '

# Row 25: 'Aptos' -> 'PEPE'
$ws.Range('B25').Value = 'PEPE'
$ws.Range('C25').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D25').Value = '0.0₃0808'
$ws.Range('E25').Value = '  +1.96%  '

# Row 26: 'PEPE' -> 'InternetComputer(DFINITY)'
$ws.Range('B26').Value = 'InternetComputer(DFINITY)'
$ws.Range('C26').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$cell = $ws.Range('D26')
$cell.NumberFormat = "@"
$cell.Value = '7.09'
$cell.Style = "Normal"
$ws.Range('E26').Value = '  +1.17%  '

# Row 27: 'InternetComputer(DFINITY)' -> 'FirstDigitalUSD'
$ws.Range('B27').Value = 'FirstDigitalUSD'
$ws.Range('C27').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$cell = $ws.Range('D27')
$cell.NumberFormat = "@"
$cell.Value = '1.00'
$cell.Style = "Normal"
$ws.Range('E27').Value = '  -0.02%  '

# Row 28: 'FirstDigitalUSD' -> 'Bittensor'
$ws.Range('B28').Value = 'Bittensor'
$ws.Range('C28').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$cell = $ws.Range('D28')
$cell.NumberFormat = "@"
$cell.Value = '420.59'
$cell.Style = "Normal"
$ws.Range('E28').Value = '  -1.88%  '

# Row 29: 'Bittensor' -> 'Fetch.AI'
$ws.Range('B29').Value = 'Fetch.AI'
$ws.Range('C29').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$cell = $ws.Range('D29')
$cell.NumberFormat = "@"
$cell.Value = '1.13'
$cell.Style = "Normal"
$ws.Range('E29').Value = '  +3.06%  '

# Row 30: 'Fetch.AI' -> 'PancakeSwap'
$ws.Range('B30').Value = 'PancakeSwap'
$ws.Range('C30').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$cell = $ws.Range('D30')
$cell.NumberFormat = "@"
$cell.Value = '1.59'
$cell.Style = "Normal"
$ws.Range('E30').Value = '  +0.41%  '

# Row 31: 'PancakeSwap' -> 'Monero'
$ws.Range('B31').Value = 'Monero'
$ws.Range('C31').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$cell = $ws.Range('D31')
$cell.NumberFormat = "@"
$cell.Value = '160.31'
$cell.Style = "Normal"
$ws.Range('E31').Value = '  +2.32%  '

# Row 32: 'Monero' -> 'WhiteBITCoin'
$ws.Range('B32').Value = 'WhiteBITCoin'
$ws.Range('C32').Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$cell = $ws.Range('D32')
$cell.NumberFormat = "@"
$cell.Value = '18.93'
$cell.Style = "Normal"
$ws.Range('E32').Value = '  -0.33%  '

# Row 33: 'WhiteBITCoin' -> 'USDe'
$ws.Range('B33').Value = 'USDe'
$ws.Range('C33').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$cell = $ws.Range('D33')
$cell.NumberFormat = "@"
$cell.Value = '1.00'
$cell.Style = "Normal"
$ws.Range('E33').Value = '  +0.05%  '

# Row 34: 'USDe' -> 'EthereumClassic'
$ws.Range('B34').Value = 'EthereumClassic'
$ws.Range('C34').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$cell = $ws.Range('D34')
$cell.NumberFormat = "@"
$cell.Value = '17.79'
$cell.Style = "Normal"
$ws.Range('E34').Value = '  +1.21%  '

# Row 35: 'EthereumClassic' -> 'Kaspa'
$ws.Range('B35').Value = 'Kaspa'
$ws.Range('C35').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$cell = $ws.Range('D35')
$cell.NumberFormat = "@"
$cell.Value = '0.103'
$cell.Style = "Normal"
$ws.Range('E35').Value = '  -3.94%  '

# Row 36: 'Kaspa' -> 'PolygonEcosystemToken'
$ws.Range('B36').Value = 'PolygonEcosystemToken'
$ws.Range('C36').Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$cell = $ws.Range('D36')
$cell.NumberFormat = "@"
$cell.Value = '0.293'
$cell.Style = "Normal"
$ws.Range('E36').Value = '  -1.19%  '

# Row 37: 'PolygonEcosystemToken' -> 'RenderToken'
$ws.Range('B37').Value = 'RenderToken'
$ws.Range('C37').Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$cell = $ws.Range('D37')
$cell.NumberFormat = "@"
$cell.Value = '4.26'
$cell.Style = "Normal"
$ws.Range('E37').Value = '  -1.50%  '

# Row 38: 'RenderToken' -> 'Stacks'
$ws.Range('B38').Value = 'Stacks'
$ws.Range('C38').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$cell = $ws.Range('D38')
$cell.NumberFormat = "@"
$cell.Value = '1.46'
$cell.Style = "Normal"
$ws.Range('E38').Value = '  +3.01%  '

# Row 39: 'Stacks' -> 'ImmutableX'
$ws.Range('B39').Value = 'ImmutableX'
$ws.Range('C39').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$cell = $ws.Range('D39')
$cell.NumberFormat = "@"
$cell.Value = '1.07'
$cell.Style = "Normal"
$ws.Range('E39').Value = '  +0.19%  '

# Row 40: 'ImmutableX' -> 'dogwifhat'
$ws.Range('B40').Value = 'dogwifhat'
$ws.Range('C40').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$cell = $ws.Range('D40')
$cell.NumberFormat = "@"
$cell.Value = '2.00'
$cell.Style = "Normal"
$ws.Range('E40').Value = '  +1.44%  '

# Row 41: 'dogwifhat' -> 'Filecoin'
$ws.Range('B41').Value = 'Filecoin'
$ws.Range('C41').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$cell = $ws.Range('D41')
$cell.NumberFormat = "@"
$cell.Value = '3.33'
$cell.Style = "Normal"
$ws.Range('E41').Value = '  +1.55%  '

# Row 42: 'Filecoin' -> 'Aave'
$ws.Range('B42').Value = 'Aave'
$ws.Range('C42').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$cell = $ws.Range('D42')
$cell.NumberFormat = "@"
$cell.Value = '128.84'
$cell.Style = "Normal"
$ws.Range('E42').Value = '  -0.39%  '

# Row 43: 'Aave' -> 'Cronos'
$ws.Range('B43').Value = 'Cronos'
$ws.Range('C43').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$cell = $ws.Range('D43')
$cell.NumberFormat = "@"
$cell.Value = '0.0710'
$cell.Style = "Normal"
$ws.Range('E43').Value = '  +0.15%  '

# Row 44: 'Cronos' -> 'ARBITRUM'
$ws.Range('B44').Value = 'ARBITRUM'
$ws.Range('C44').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$cell = $ws.Range('D44')
$cell.NumberFormat = "@"
$cell.Value = '0.478'
$cell.Style = "Normal"
$ws.Range('E44').Value = '  +0.91%  '

# Row 45: 'ARBITRUM' -> 'Mantle'
$ws.Range('B45').Value = 'Mantle'
$ws.Range('C45').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$cell = $ws.Range('D45')
$cell.NumberFormat = "@"
$cell.Value = '0.554'
$cell.Style = "Normal"
$ws.Range('E45').Value = '  +0.74%  '

# Row 46: 'Mantle' -> 'Stellar'
$ws.Range('B46').Value = 'Stellar'
$ws.Range('C46').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$cell = $ws.Range('D46')
$cell.NumberFormat = "@"
$cell.Value = '0.0911'
$cell.Style = "Normal"
$ws.Range('E46').Value = '  +1.56%  '

# Row 47: 'Stellar' -> 'BitgetToken'
$ws.Range('B47').Value = 'BitgetToken'
$ws.Range('C47').Value = 'https://coinranking.com/coin/q7gMmMdLb+bitgettoken-bgb'
$cell = $ws.Range('D47')
$cell.NumberFormat = "@"
$cell.Value = '1.11'
$cell.Style = "Normal"
$ws.Range('E47').Value = '  +0.56%  '

# Row 48: 'BitgetToken' -> 'Optimism'
$ws.Range('B48').Value = 'Optimism'
$ws.Range('C48').Value = 'https://coinranking.com/coin/n1p-s_gm1+optimism-op'
$cell = $ws.Range('D48')
$cell.NumberFormat = "@"
$cell.Value = '1.34'
$cell.Style = "Normal"
$ws.Range('E48').Value = '  -4.62%  '

# Row 49: 'Optimism' -> 'InjectiveProtocol'
$ws.Range('B49').Value = 'InjectiveProtocol'
$ws.Range('C49').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$cell = $ws.Range('D49')
$cell.NumberFormat = "@"
$cell.Value = '16.58'
$cell.Style = "Normal"
$ws.Range('E49').Value = '  +0.49%  '

# Row 50: 'InjectiveProtocol' -> 'BabyDogeCoin'
$ws.Range('B50').Value = 'BabyDogeCoin'
$ws.Range('C50').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D50').Value = '0.0₆0204'
$ws.Range('E50').Value = '  +5.84%  '

# Row 51: 'BabyDogeCoin' -> 'Hedera'
$ws.Range('B51').Value = 'Hedera'
$ws.Range('C51').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$cell = $ws.Range('D51')
$cell.NumberFormat = "@"
$cell.Value = '0.0428'
$cell.Style = "Normal"
$ws.Range('E51').Value = '  +1.88%  '
